# FBKeyWordDrivenLogin.xlsx
#
# "Added: Functionality to check the KeyWordDriven FrameWork from the
#  facebook website"
#
# On the "Login" sheet, the "Enter Password" row (row 5) used to locate
# the password field by id="pass". Switch it to a CSS selector locator
# (#pass) instead:
#   Locator (col B)      : "id "  -> "cssSelector"
#   LocatorValue (col C) : "pass" -> "#pass"
#
# Also restore the workbook's on-screen focus to the Login sheet /
# cell B5 (the row that was just edited), matching the saved view
# state in the workbook.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# Set LocatorValue before Locator so new shared strings land in the
# same order ("#pass" then "cssSelector") as the authored workbook.
$loginSheet.Range("C5").Value = "#pass"
$loginSheet.Range("B5").Value = "cssSelector"

# Bring the Login sheet to the front and leave B5 selected.
$loginSheet.Activate() | Out-Null
$loginSheet.Range("B5").Select() | Out-Null
